# Append a new listing as row 9 on the "ランサーズ" sheet, pushing the
# previous row 9 ("UTAGE" listing) down to row 10, refresh every "取得日時"
# timestamp, and bump row 7's transaction-period text by one day.
# Matches commit message: "Append: 2026-01-29 01:51 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-29 01:51:02"

# ---------------------------------------------------------------------
# 1) Push the existing row 9 (UTAGE listing) down to row 10 first, so it
#    is not lost once row 9 gets overwritten with the new listing below.
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = $newTimestamp
$ws.Cells.Item(10, 2).Value = $ws.Cells.Item(9, 2).Value2
$ws.Cells.Item(10, 3).Value = $ws.Cells.Item(9, 3).Value2
$ws.Cells.Item(10, 4).Value = $ws.Cells.Item(9, 4).Value2
$ws.Cells.Item(10, 5).Value = $ws.Cells.Item(9, 5).Value2
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5480853"
$ws.Cells.Item(10, 7).Value = $ws.Cells.Item(9, 7).Value2

# ---------------------------------------------------------------------
# 2) Overwrite row 9 with the brand-new listing.
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = $newTimestamp
$ws.Cells.Item(9, 2).Value = "【短期/読み取り専用】既存システムに触れない行動ログ基盤の構築(Fintech系)"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5481091"
$ws.Cells.Item(9, 7).Value = 33

# ---------------------------------------------------------------------
# 3) Refresh the "取得日時" timestamp for every other pre-existing row.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# ---------------------------------------------------------------------
# 4) Row 7's price cell gained one extra transaction-period day.
# ---------------------------------------------------------------------
$ws.Cells.Item(7, 4).Value = "10,000 円 ~ 20,000 円 / 募集期間 3 日、取引期間 1 日"

# ---------------------------------------------------------------------
# 5) Rebuild the F-column hyperlinks (F2:F10) from scratch so every link
#    points at the right URL, then restyle the cells back to the
#    workbook's built-in "Hyperlink" look.
# ---------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$urls = @{
    2 = "https://www.lancers.jp/work/detail/5455098"
    3 = "https://www.lancers.jp/work/detail/5480502"
    4 = "https://www.lancers.jp/work/detail/5445159"
    5 = "https://www.lancers.jp/work/detail/5445154"
    6 = "https://www.lancers.jp/work/detail/5480464"
    7 = "https://www.lancers.jp/work/detail/5480298"
    8 = "https://www.lancers.jp/work/detail/5480933"
    9 = "https://www.lancers.jp/work/detail/5481091"
    10 = "https://www.lancers.jp/work/detail/5480853"
}

foreach ($r in 2..10) {
    $target = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($target, $urls[$r]) | Out-Null
    $target.Style = "Hyperlink"
}
